# Corrección en la lógica de compra de simultaneas
# Swap the values of "cantidad" (E) and "monto" (F) columns for the
# SIMULTANEA purchase rows (rows 2-9), since they had been entered
# in the wrong columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $eCell = $ws.Cells.Item($row, 5)  # column E - cantidad
    $fCell = $ws.Cells.Item($row, 6)  # column F - monto

    $eValue = $eCell.Value()
    $fValue = $fCell.Value()

    $eCell.Value = $fValue
    $fCell.Value = $eValue
}
